# Tripadvisor New Orleans shard 126 - update:
#  1) Reorder worksheets so "review_info" is the first tab and "hotel_info"
#     is the second tab (previously hotel_info was first).
#  2) Add a new "State" column to hotel_info, inserted between Hotel_Name
#     and City, with value "Louisiana" for the existing data row.

$wb = $excel.ActiveWorkbook

# --- 1) Reorder the sheets -------------------------------------------------
$reviewInfo = $wb.Worksheets.Item("review_info")
$hotelInfo  = $wb.Worksheets.Item("hotel_info")
$reviewInfo.Move($hotelInfo)

# NOTE: worksheet object handles in this host are positional, so after the
# Move() above the old $hotelInfo / $reviewInfo variables now resolve to
# whatever sheet occupies that slot. Re-fetch by name before continuing.

# --- 2) Insert the "State" column into hotel_info --------------------------
$hotelInfo = $wb.Worksheets.Item("hotel_info")
$hotelInfo.Columns.Item(3).Insert()
$hotelInfo.Cells.Item(1, 3).Value = "State"
$hotelInfo.Cells.Item(2, 3).Value = "Louisiana"
